# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row (in column F) -> new value mapping, applied identically to both sheets.
$updates = @{
    3  = 46
    5  = 82
    7  = 1224
    8  = 1512
    10 = 376
    12 = 141
    17 = 293
    18 = 319
    19 = 1710
    20 = 64
    22 = 173
    25 = 330
    26 = 4118
    30 = 1072
    31 = 132
    33 = 481
    35 = 218
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
